# Update the four controller transfer-function strings in row 2.
# (PowerShell single-quoted strings are used so the literal "$" characters
# in the LaTeX-style formulas are NOT treated as variable interpolation.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$C(s) =\frac{3.174\,s^2+16.08\,s+13.18}{0.0023\,s^2+1.219\,s}$'
$ws.Range("C2").Value = '$C_1(s) =\frac{3.631}{s}$  $C_2(s) =\frac{2.826\,s+6.76}{0.02427\,s+1.0}$'
$ws.Range("D2").Value = '$C_1(s) =4.4585$  $C_2(s) =\frac{0.6278}{s}$  $C_3(s) =\frac{1.823\,s}{0.01132\,s+1.0}$'
$ws.Range("E2").Value = '$C(s) =\frac{0.002797\,s^4+1.3\,s^3+20.24\,s^2+83.19\,s+64.24}{8.676e-11\,s^4+1.025e-6\,s^3+0.00325\,s^2+1.293\,s}$'

# Update the numeric data table (B3:E20) with the refreshed benchmark numbers.
$ws.Range("B3").Value = 0.10109578804696412
$ws.Range("C3").Value = 2.1374021739348059
$ws.Range("D3").Value = 1.5929736852202085
$ws.Range("E3").Value = 0.020838485521172914

$ws.Range("B4").Value = 13.175118108275534
$ws.Range("C4").Value = 6.7601334600593184
$ws.Range("D4").Value = 4.4585431936982811
$ws.Range("E4").Value = 64.238932755179121

$ws.Range("B5").Value = 1.218692506687403
$ws.Range("C5").Value = 1.861924929648904
$ws.Range("D5").Value = 7.1021397489810179
$ws.Range("E5").Value = 1.2925635821932677

$ws.Range("B6").Value = 0.195780770592273
$ws.Range("C6").Value = 0.39370388968740699
$ws.Range("D6").Value = 0.40893142255738102
$ws.Range("E6").Value = 0.24129851120039555

$ws.Range("E7").Value = 0.015574343636600292

$ws.Range("B10").Value = 103.7209327665792
$ws.Range("C10").Value = 16.222375069598719
$ws.Range("D10").Value = 36.129690930396485
$ws.Range("E10").Value = 111.58996146053157

$ws.Range("E11").Value = 88.397527344577796

$ws.Range("B12").Value = 1.3999995754524437
$ws.Range("C12").Value = 1.1590725796863817
$ws.Range("D12").Value = 1.1165562229787664
$ws.Range("E12").Value = 1.3999899061207535

$ws.Range("B13").Value = 0.11301156132123021
$ws.Range("C13").Value = 3.60140564251454
$ws.Range("D13").Value = 2.5637228096197058
$ws.Range("E13").Value = 0.015122524846155751

$ws.Range("B14").Value = 0.48068989011372149
$ws.Range("C14").Value = 6.4460571862463647
$ws.Range("D14").Value = 15.310819017254429
$ws.Range("E14").Value = 0.19386120511433746

$ws.Range("B15").Value = 4.2640505337536716
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 3.6823269017976612

$ws.Range("B17").Value = 1.0426405053375365
$ws.Range("C17").Value = 0.99912880337904664
$ws.Range("D17").Value = 0.99942353561422681
$ws.Range("E17").Value = 1.0368232690179766

$ws.Range("B18").Value = 0.23033952318255399
$ws.Range("C18").Value = 10.802845232701818
$ws.Range("D18").Value = 44.025881970460809
$ws.Range("E18").Value = 0.029604821498780413

$ws.Range("B19").Value = 9.4641072432786117
$ws.Range("C19").Value = 38.295282709793206
$ws.Range("D19").Value = 14.402026366559678
$ws.Range("E19").Value = 28.041576530166026

$ws.Range("B20").Value = 62.928504501873149
$ws.Range("C20").Value = 76.944318229849955
$ws.Range("D20").Value = 90.831736111463073
$ws.Range("E20").Value = 60.213184619523339

# Re-intern the header/label styling so the previously-duplicated border/
# number-format style slots collapse into the first matching slot (style
# index 1 instead of 3), matching the canonical XML's compacted style table.
$ws.Range("A1:E2").NumberFormat = "@"
$ws.Range("A3:A20").NumberFormat = "@"

# Restore the small column-width tweaks from the commit (ColumnWidth is
# expressed in character units and Excel rounds it to whole pixels
# internally, so the inputs below are chosen to land on the closest
# achievable stored width to the canonical 55 / 63 / 74.7109375 / 103.140625).
$ws.Columns.Item(2).ColumnWidth = 54.166666666666664
$ws.Columns.Item(3).ColumnWidth = 62.166666666666664
$ws.Columns.Item(4).ColumnWidth = 73.83333333333334
$ws.Columns.Item(5).ColumnWidth = 102.33333333333334
